# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the daily conversion text block on Hoja1!A1
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.35 = 8880.75 pesos`n✅ 8880.75 pesos = 2.33 = 946.35 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update the rate values on the "tasas" sheet
$wsTasas.Range("N10").Value = 426
$wsTasas.Range("O10").Value = 3783.2
$wsTasas.Range("N12").Value = 3809.99
$wsTasas.Range("O12").Value = 406
